$wb = $excel.ActiveWorkbook

# --- Update the daily conversion text on "Hoja1" ---
$ws1 = $wb.Worksheets.Item("Hoja1")
$text = $ws1.Range("A1").Value()
$text = $text -replace [regex]::Escape("1000 Bs = 1.61 = 5935.1 pesos"), "1000 Bs = 1.62 = 6043.95 pesos"
$text = $text -replace [regex]::Escape("5935.1 pesos = 1.6 = 962.97 Bs"), "6043.95 pesos = 1.62 = 969.16 Bs"
$ws1.Range("A1").Value = $text

# --- Update the tasa values on "tasas" ---
$ws2 = $wb.Worksheets.Item("tasas")
$ws2.Range("N10").Value = 618
$ws2.Range("O10").Value = 3735.16
$ws2.Range("N12").Value = 3730
$ws2.Range("O12").Value = 598.1130000000001
